$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E, shifting Writing Tip/Writing Tip
# Illustration/Reading Tip/Tags one column to the right, and use it for a
# new "Romaji Extra" field.
$ws.Columns("E").Insert()
$ws.Columns("E").ColumnWidth = 22.42578125

$ws.Range("E1").Value = "Romaji Extra"
$ws.Range("E53").Value = "<small><small>/ztsa/</small></small>"
$ws.Range("E160").Value = "<small><small>/ztsa/</small></small>"
